$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "ファー・ローズトゥ・ロードリプレイ-RPGセッションガイド"
$ws.Range("C17").Value = "Far Rose to Road Replay RPG Session Guide"
$ws.Range("D17").Value = "Softbank Books"
$ws.Range("E17").Value = "session_guide.jpg"
$ws.Range("F17").Value = "replay"

$ws.Range("A18").Value = 2002
$ws.Range("B18").Value = "ローズ・トゥ・ロード"
$ws.Range("C18").Value = "Roads to Lord"
$ws.Range("D18").Value = "Enterrain"
$ws.Range("E18").Value = "roads_to_lord_2002.jpg"
$ws.Range("F18").Value = "rulebook"

$ws.Range("B19").Value = "ザ・ストレンジソング"
$ws.Range("C19").Value = "The Strange Song"
$ws.Range("D19").Value = "Arclight"
$ws.Range("E19").Value = "strange_song.jpg"
$ws.Range("F19").Value = "supplement"

$ws.Range("B20").Value = "タトゥーノ〜“風に絵を書く”かりそめの魔法〜"
$ws.Range("C20").Value = "Tatuno ~ `"Drawing a picture in the wind`" Karisome's magic ~"
$ws.Range("D20").Value = "Arclight"
$ws.Range("E20").Value = "tatuno.jpg"
$ws.Range("F20").Value = "supplement"

$ws.Range("B21").Value = "ゲームマスター・スクリーン〜忘却の呪縛、近づく頃〜"
$ws.Range("C21").Value = "Gamemaster's Screen-The Curse of Oblivion, When It's Approaching-"
$ws.Range("D21").Value = "Arclight"
$ws.Range("E21").Value = "gamemasters_screen.jpg"
$ws.Range("F21").Value = "supplement"

$ws.Range("A22").Value = 2006
$ws.Range("B22").Value = "ソングシーカー"
$ws.Range("C22").Value = "Song Seeker"
$ws.Range("D22").Value = "Shinkigensha"
$ws.Range("E22").Value = "song_seeker.jpg"
$ws.Range("F22").Value = "replay"

$ws.Range("A23").Value = 2010
$ws.Range("B23").Value = "ローズ・トゥ・ロード"
$ws.Range("C23").Value = "Roads to Lord"
$ws.Range("D23").Value = "Enterrain"
$ws.Range("E23").Value = "roads_to_lord_2010.jpg"
$ws.Range("F23").Value = "rulebook"

$ws.Rows.Item(24).Select()